{"js": "// Change the answer to \"How many times have you submitted this assessment\n// (including this time)?\" from \"Once\" to \"Twice\", keeping the existing bold\n// formatting of the answer intact.\nconst body = context.document.body;\n\nconst results = body.search(\"Once\", { matchCase: true, matchWholeWord: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find the text \"Once\" to replace.');\n}\n\n// Replace in place so the run's formatting (bold) is preserved.\nresults.items[0].insertText(\"Twice\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Change the answer to \"How many times have you submitted this assessment\n# (including this time)?\" from \"Once\" to \"Twice\", keeping the existing bold\n# formatting of the answer intact.\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"Once\"\n$find.MatchWholeWord = $true\n$find.MatchCase = $true\n\nif ($find.Execute()) {\n    # Assigning .Text on the found range replaces the text in place,\n    # preserving the run's existing character formatting (bold).\n    $range.Text = \"Twice\"\n} else {\n    throw 'Could not find the text \"Once\" to replace.'\n}\n"}
